# Insert a new weekly price row for "Feria Lagunitas de Puerto Montt - Zapallo"
# above the current row 229, shifting all subsequent rows down by one
# (old row 229 becomes new row 230, ..., old row 316 becomes new row 317).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("229:229").Insert()

$ws.Range("A229").Value = 4
$ws.Range("B229").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C229").Value = "Los Lagos"
$ws.Range("D229").Value = 44726
$ws.Range("E229").Value = 10
$ws.Range("F229").Value = 100112045
$ws.Range("G229").Value = "Zapallo"
$ws.Range("H229").Value = "Paine"
$ws.Range("I229").Value = "1a (guarda)"
$ws.Range("J229").Value = 1000
$ws.Range("K229").Value = 500
$ws.Range("L229").Value = 500
$ws.Range("M229").Value = 500
$ws.Range("N229").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O229").Value = "Región de O'Higgins"
$ws.Range("P229").Value = 500
$ws.Range("Q229").Value = 1
$ws.Range("R229").Value = "Hortaliza"
